$d = $word.ActiveDocument

function Set-CellText($cell, [string]$newText) {
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# --- Table 1: "Установленные права и разрешённые действия для групп" ---
# Rows 3..9 hold the permission-bit values repeated in columns 1 and 2
# (e.g. "100","100" -> "010","010"), leaving every other cell untouched.
$table1 = $d.Tables.Item(1)

$table1Map = @{
    3 = "010"
    4 = "020"
    5 = "030"
    6 = "040"
    7 = "050"
    8 = "060"
    9 = "070"
}

foreach ($rowIndex in $table1Map.Keys) {
    $newValue = $table1Map[$rowIndex]
    Set-CellText $table1.Cell($rowIndex, 1) $newValue
    Set-CellText $table1.Cell($rowIndex, 2) $newValue
}

# --- Table 2: "Минимальные права для совершения операций" ---
# Rows 2..8 hold "Минимальные права на директорию" (col 2) and
# "Минимальные права на файл" (col 3) values.
$table2 = $d.Tables.Item(2)

$table2Map = @{
    2 = @{ Dir = "030"; File = "020" }  # Создание файла
    3 = @{ Dir = "030"; File = "020" }  # Удаление файла
    4 = @{ Dir = "050"; File = "040" }  # Чтение файла
    5 = @{ Dir = "030"; File = "020" }  # Запись в файл
    6 = @{ Dir = "030"; File = "020" }  # Переименование файла
    7 = @{ Dir = "030"; File = "030" }  # Создание поддиректории
    8 = @{ Dir = "030"; File = "030" }  # Удаление поддиректории
}

foreach ($rowIndex in $table2Map.Keys) {
    $values = $table2Map[$rowIndex]
    Set-CellText $table2.Cell($rowIndex, 2) $values.Dir
    Set-CellText $table2.Cell($rowIndex, 3) $values.File
}

Write-Output "done"
